$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shared string rename: "City" -> "Location" ---
$ws.Range("A2").Value = "Location"

# --- 2. Row 3 (Quality) - new cells, copy style from row 4 (style 11) ---
$ws.Range("B4").Copy($ws.Range("B3"))
$ws.Range("C4").Copy($ws.Range("C3"))
$ws.Range("E4").Copy($ws.Range("E3"))
$ws.Range("F4").Copy($ws.Range("F3"))
$ws.Range("H4").Copy($ws.Range("H3"))
$ws.Range("I4").Copy($ws.Range("I3"))
$ws.Range("K4").Copy($ws.Range("K3"))
$ws.Range("L4").Copy($ws.Range("L3"))

# --- 3. Row 5 (Chance of Winning) - new cells, copy style from row 2 (style 10) ---
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("C2").Copy($ws.Range("C5"))
$ws.Range("E2").Copy($ws.Range("E5"))
$ws.Range("F2").Copy($ws.Range("F5"))
$ws.Range("H2").Copy($ws.Range("H5"))
$ws.Range("I2").Copy($ws.Range("I5"))
$ws.Range("K2").Copy($ws.Range("K5"))
$ws.Range("L2").Copy($ws.Range("L5"))

# --- 4. Set values for rows 3-26 per the target data ---
# Row 3
$ws.Range("B3").Value = 0.263291367343765
$ws.Range("E3").Value = 0.5057489238611556
$ws.Range("H3").Value = 0.4137738998947945
$ws.Range("K3").Value = 0.7627585613100407

# Row 4
$ws.Range("B4").Value = 0.9131806738518395
$ws.Range("E4").Value = 0.9883674465155902
$ws.Range("H4").Value = 0.9906281334855674
$ws.Range("K4").Value = 0.9995591623406129

# Row 5
$ws.Range("B5").Value = 24.04325882503515
$ws.Range("E5").Value = 49.9865772454658
$ws.Range("H5").Value = 40.98960661378243
$ws.Range("K5").Value = 76.24223086111954

# Row 6
$ws.Range("B6").Value = 0.3283006
$ws.Range("C6").Value = 0.6716993999999999
$ws.Range("E6").Value = 0.4365911
$ws.Range("F6").Value = 0.5634089
$ws.Range("H6").Value = 0.4430702
$ws.Range("I6").Value = 0.5569298
$ws.Range("K6").Value = 0.4876401
$ws.Range("L6").Value = 0.5123599

# Row 7
$ws.Range("B7").Value = 14.9962966
$ws.Range("C7").Value = 21.4901812
$ws.Range("E7").Value = 23.991193
$ws.Range("F7").Value = 26.994349
$ws.Range("H7").Value = 28.8108116
$ws.Range("I7").Value = 31.5167234
$ws.Range("K7").Value = 24.9978058
$ws.Range("L7").Value = 25.5084802

# Row 8
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 6
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 7
$ws.Range("H8").Value = 8
$ws.Range("I8").Value = 10
$ws.Range("K8").Value = 7
$ws.Range("L8").Value = 7

# Row 9
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 9
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 10
$ws.Range("H9").Value = 11
$ws.Range("I9").Value = 14
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = 10

# Row 10
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 10
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 14
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = 17
$ws.Range("K10").Value = 13
$ws.Range("L10").Value = 13

# Row 11
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 13
$ws.Range("E11").Value = 14
$ws.Range("F11").Value = 16
$ws.Range("H11").Value = 16
$ws.Range("I11").Value = 20
$ws.Range("K11").Value = 14
$ws.Range("L11").Value = 14

# Row 12
$ws.Range("B12").Value = 7
$ws.Range("C12").Value = 13
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 17
$ws.Range("H12").Value = 18
$ws.Range("I12").Value = 21
$ws.Range("K12").Value = 16
$ws.Range("L12").Value = 17

# Row 13
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 16
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = 20
$ws.Range("H13").Value = 21
$ws.Range("I13").Value = 24
$ws.Range("K13").Value = 17
$ws.Range("L13").Value = 17

# Row 14
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 16
$ws.Range("E14").Value = 17
$ws.Range("F14").Value = 21
$ws.Range("H14").Value = 22
$ws.Range("I14").Value = 25
$ws.Range("K14").Value = 19
$ws.Range("L14").Value = 20

# Row 15
$ws.Range("B15").Value = 12
$ws.Range("C15").Value = 17
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 23
$ws.Range("H15").Value = 24
$ws.Range("I15").Value = 27
$ws.Range("K15").Value = 21
$ws.Range("L15").Value = 21

# Row 16
$ws.Range("B16").Value = 13
$ws.Range("C16").Value = 19
$ws.Range("E16").Value = 21
$ws.Range("F16").Value = 24
$ws.Range("H16").Value = 26
$ws.Range("I16").Value = 29
$ws.Range("K16").Value = 22
$ws.Range("L16").Value = 23

# Row 17
$ws.Range("B17").Value = 14
$ws.Range("C17").Value = 20
$ws.Range("E17").Value = 23
$ws.Range("F17").Value = 26
$ws.Range("H17").Value = 28
$ws.Range("I17").Value = 30
$ws.Range("K17").Value = 24
$ws.Range("L17").Value = 24

# Row 18
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = 22
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 27
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 32
$ws.Range("K18").Value = 26
$ws.Range("L18").Value = 27

# Row 19
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = 23
$ws.Range("E19").Value = 27
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = 31
$ws.Range("I19").Value = 34
$ws.Range("K19").Value = 27
$ws.Range("L19").Value = 28

# Row 20
$ws.Range("B20").Value = 17
$ws.Range("C20").Value = 24
$ws.Range("E20").Value = 28
$ws.Range("F20").Value = 31
$ws.Range("H20").Value = 33
$ws.Range("I20").Value = 36
$ws.Range("K20").Value = 29
$ws.Range("L20").Value = 30

# Row 21
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 33
$ws.Range("H21").Value = 35
$ws.Range("I21").Value = 38
$ws.Range("K21").Value = 31
$ws.Range("L21").Value = 31

# Row 22
$ws.Range("B22").Value = 20
$ws.Range("C22").Value = 28
$ws.Range("E22").Value = 31
$ws.Range("F22").Value = 35
$ws.Range("H22").Value = 38
$ws.Range("I22").Value = 40
$ws.Range("K22").Value = 33
$ws.Range("L22").Value = 34

# Row 23
$ws.Range("B23").Value = 23
$ws.Range("C23").Value = 30
$ws.Range("E23").Value = 34
$ws.Range("F23").Value = 37
$ws.Range("H23").Value = 40
$ws.Range("I23").Value = 43
$ws.Range("K23").Value = 35
$ws.Range("L23").Value = 35

# Row 24
$ws.Range("B24").Value = 24
$ws.Range("C24").Value = 33
$ws.Range("E24").Value = 38
$ws.Range("F24").Value = 41
$ws.Range("H24").Value = 44
$ws.Range("I24").Value = 46
$ws.Range("K24").Value = 38
$ws.Range("L24").Value = 38

# Row 25
$ws.Range("B25").Value = 27
$ws.Range("C25").Value = 36
$ws.Range("E25").Value = 41
$ws.Range("F25").Value = 44
$ws.Range("H25").Value = 48
$ws.Range("I25").Value = 50
$ws.Range("K25").Value = 42
$ws.Range("L25").Value = 42

# Row 26
$ws.Range("B26").Value = 31
$ws.Range("C26").Value = 41
$ws.Range("E26").Value = 45
$ws.Range("F26").Value = 50
$ws.Range("H26").Value = 54
$ws.Range("I26").Value = 56
$ws.Range("K26").Value = 47
$ws.Range("L26").Value = 48

# --- 5. Merge new cells B3:C3, E3:F3, H3:I3, K3:L3, B5:C5, E5:F5, H5:I5, K5:L5 ---
$ws.Range("B3:C3").Merge()
$ws.Range("E3:F3").Merge()
$ws.Range("H3:I3").Merge()
$ws.Range("K3:L3").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("E5:F5").Merge()
$ws.Range("H5:I5").Merge()
$ws.Range("K5:L5").Merge()
